$wb = $excel.ActiveWorkbook

$wsLab    = $wb.Worksheets.Item("Translations - Lab")
$wsMarket = $wb.Worksheets.Item("Translations - Market")

# --- Relocate the two "liquid" translations ("lab.market.menu" and
#     "lab.liquid.index.title") from "Translations - Market" (rows 68-69)
#     into "Translations - Lab" (new rows 25-26), and append the brand new
#     "liquid creation" strings right after them (rows 27-32). The feature
#     is now backed by transactions from the pricelist, hence the new keys.

# Give the eight new rows the same formatting as the last existing row.
$fmtSource = $wsLab.Range("A24:C24")
$newRange  = $wsLab.Range("A25:C32")
$fmtSource.Copy()
$newRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$newRows = @(
    @("cs", "lab.market.menu",           "Tržiště"),
    @("cs", "lab.liquid.index.title",    "Liquidy"),
    @("cs", "lab.liquid.create.success", "Liquid [{{data.name}}] byl úspěšně vytvořen za {{data.transaction.amount}} puffíků."),
    @("cs", "lab.liquid.create.title",   "Nový liquid"),
    @("cs", "lab.liquid.name.label",     "Název"),
    @("cs", "lab.liquid.volume.label",   "Celkový objem"),
    @("cs", "lab.liquid.mixed.label",    "Datum míchání"),
    @("cs", "lab.liquid.create",         "Uložit")
)

$startRow = 25
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $wsLab.Range("A$r").Value = $newRows[$i][0]
    $wsLab.Range("B$r").Value = $newRows[$i][1]
    $wsLab.Range("C$r").Value = $newRows[$i][2]
}

# Drop the two rows that were moved away from the Market sheet.
$wsMarket.Rows("68:69").Delete()

# Market used to be the selected/active sheet (C58 selected); touch its
# selection first so the change is recorded there too.
$wsMarket.Activate()
$wsMarket.Range("B58").Select() | Out-Null

# Lab is now the active/selected sheet, with B27 (the first brand-new row)
# selected.
$wsLab.Activate()
$wsLab.Range("B27").Select() | Out-Null
